$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking price cells so Excel keeps them as
# literal text (matching the original "Price" column formatting) instead of
# auto-converting to numbers on assignment.
$textCells = @("D5","D6","D7","D9","D10","D11","D12","D13","D14","D19","D20","D21","D22","D24","D26","D28","D29","D30","D31","D32","D33","D35","D36","D37","D38","D40","D41","D43","D44","D47","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '42.397.99'
$ws.Range('E2').Value = '  +0.69%  '
$ws.Range('D3').Value = '2.229.19'
$ws.Range('E3').Value = '  -0.17%  '
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').Value = '245.07'
$ws.Range('E5').Value = '  -0.11%  '
$ws.Range('D6').Value = '0.629'
$ws.Range('E6').Value = '  +1.50%  '
$ws.Range('D7').Value = '74.20'
$ws.Range('E7').Value = '  -1.63%  '
$ws.Range('E8').Value = '  +0.11%  '
$ws.Range('D9').Value = '0.617'
$ws.Range('E9').Value = '  -0.06%  '
$ws.Range('D10').Value = '42.98'
$ws.Range('E10').Value = '  +4.09%  '
$ws.Range('D11').Value = '0.0963'
$ws.Range('E11').Value = '  +2.32%  '
$ws.Range('D12').Value = '7.15'
$ws.Range('E12').Value = '  +1.02%  '
$ws.Range('D13').Value = '0.104'
$ws.Range('E13').Value = '  +0.76%  '
$ws.Range('D14').Value = '14.47'
$ws.Range('E14').Value = '  -0.21%  '
$ws.Range('E15').Value = '  +0.43%  '
$ws.Range('D16').Value = '2.237.48'
$ws.Range('E16').Value = '  +0.32%  '
$ws.Range('D17').Value = '42.166.57'
$ws.Range('E17').Value = '  +0.61%  '
$ws.Range('E18').Value = '  +13.43%  '
$ws.Range('D19').Value = '6.16'
$ws.Range('E19').Value = '  +2.13%  '
$ws.Range('D20').Value = '72.17'
$ws.Range('E20').Value = '  +1.14%  '
$ws.Range('D21').Value = '10.01'
$ws.Range('E21').Value = '  +37.72%  '
$ws.Range('D22').Value = '231.59'
$ws.Range('E22').Value = '  +0.97%  '
$ws.Range('E23').Value = '  -4.44%  '
$ws.Range('D24').Value = '11.84'
$ws.Range('E24').Value = '  +6.28%  '
$ws.Range('E25').Value = '  +0.11%  '
$ws.Range('D26').Value = '3.64'
$ws.Range('E26').Value = '  -1.04%  '
$ws.Range('E27').Value = '  +1.19%  '
$ws.Range('D28').Value = '2.23'
$ws.Range('E28').Value = '  +3.38%  '
$ws.Range('D29').Value = '167.01'
$ws.Range('E29').Value = '  -0.86%  '
$ws.Range('D30').Value = '21.05'
$ws.Range('E30').Value = '  +2.84%  '
$ws.Range('D31').Value = '5.78'
$ws.Range('E31').Value = '  +18.50%  '
$ws.Range('D32').Value = '0.0807'
$ws.Range('E32').Value = '  -1.77%  '
$ws.Range('D33').Value = '0.118'
$ws.Range('E33').Value = '  -0.27%  '
$ws.Range('E34').Value = '  +0.54%  '
$ws.Range('D35').Value = '29.57'
$ws.Range('E35').Value = '  -6.14%  '
$ws.Range('D36').Value = '4.42'
$ws.Range('E36').Value = '  +0.34%  '
$ws.Range('D37').Value = '0.0306'
$ws.Range('E37').Value = '  +2.95%  '
$ws.Range('D38').Value = '13.18'
$ws.Range('E38').Value = '  -3.73%  '
$ws.Range('E39').Value = '  +1.19%  '
$ws.Range('D40').Value = '5.63'
$ws.Range('E40').Value = '  -2.00%  '
$ws.Range('D41').Value = '63.20'
$ws.Range('E41').Value = '  +5.31%  '
$ws.Range('E42').Value = '  +0.23%  '
$ws.Range('D43').Value = '8.83'
$ws.Range('E43').Value = '  +2.59%  '
$ws.Range('D44').Value = '105.19'
$ws.Range('E44').Value = '  -6.09%  '
$ws.Range('E45').Value = '  +3.06%  '
$ws.Range('E46').Value = '  -0.24%  '
$ws.Range('D47').Value = '2.39'
$ws.Range('E47').Value = '  +6.83%  '
$ws.Range('E48').Value = '  +0.77%  '
$ws.Range('E49').Value = '  +2.09%  '
$ws.Range('D51').Value = '4.06'
$ws.Range('E51').Value = '  -0.43%  '
